$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks so stale rIds / targets don't linger once cell
# values are rewritten below (Excel does not auto-sync a hyperlink's target
# when the underlying cell's displayed text is overwritten).
$ws.Hyperlinks.Delete()

# Clear the row range we are about to fully rewrite (H column in particular
# needs a clean slate since which rows carry a 'skill summary' shifts).
$ws.Range("A2:H19").ClearContents()

# Row 2
$ws.Range("A2").Value = '2025-11-11 01:21:02'
$ws.Range("B2").Value = '【募集】習慣化+目標管理を目的としたAIネイティブなWebサービスのMVP開発'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5430365'
$ws.Range("G2").Value = 378
$ws.Range("H2").Value = '🔥AI,Ai ◆開発 ◇管理'
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5430365') | Out-Null

# Row 3
$ws.Range("A3").Value = '2025-11-11 01:21:02'
$ws.Range("B3").Value = '【急募】AIシステム構築!FirebaseとOpenAI活用の専門家募集'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5431299'
$ws.Range("G3").Value = 325
$ws.Range("H3").Value = '🔥AI,Ai'
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5431299') | Out-Null

# Row 4
$ws.Range("A4").Value = '2025-11-11 01:21:02'
$ws.Range("B4").Value = '【急募】Webシステムのエンジニア募集(API・管理画面)'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5430993'
$ws.Range("G4").Value = 228
$ws.Range("H4").Value = '🔥API ◇管理'
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5430993') | Out-Null

# Row 5
$ws.Range("A5").Value = '2025-11-11 01:21:02'
$ws.Range("B5").Value = '通話機能の安定化・不具合調査/改修(React Native × Node.js)'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5430799'
$ws.Range("G5").Value = 163
$ws.Range("H5").Value = '🔥React ◆Node.js'
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5430799') | Out-Null

# Row 6
$ws.Range("A6").Value = '2025-11-11 01:21:02'
$ws.Range("B6").Value = 'Javaプログラミング研修の演習サポート講師業務【経験不問】(再掲)'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5430954'
$ws.Range("G6").Value = 85
$ws.Range("H6").Value = '★Java'
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5430954') | Out-Null

# Row 7
$ws.Range("A7").Value = '2025-11-11 01:21:02'
$ws.Range("B7").Value = '完全在宅GASエンジニア募集/課題テストからご依頼/時給1,163円~業務フロー効率化をお任せします'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '~ 5,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5416665'
$ws.Range("G7").Value = 70
$ws.Range("H7").Value = '◆効率化'
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5416665') | Out-Null

# Row 8
$ws.Range("A8").Value = '2025-11-11 01:21:02'
$ws.Range("B8").Value = '【業務効率化】SlackとHubSpotの活用支援をお願いします'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5430436'
$ws.Range("G8").Value = 70
$ws.Range("H8").Value = '◆効率化'
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5430436') | Out-Null

# Row 9
$ws.Range("A9").Value = '2025-11-11 01:21:02'
$ws.Range("B9").Value = '【自動運転プロジェクト経験者募集】実証実験・開発を推進するプロジェクトマネージャー'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5431107'
$ws.Range("G9").Value = 68
$ws.Range("H9").Value = '◆開発'
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5431107') | Out-Null

# Row 10
$ws.Range("A10").Value = '2025-11-11 01:21:02'
$ws.Range("B10").Value = '社内利用するクローラー開発'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5431051'
$ws.Range("G10").Value = 63
$ws.Range("H10").Value = '◆開発'
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5431051') | Out-Null

# Row 11
$ws.Range("A11").Value = '2025-11-11 01:21:02'
$ws.Range("B11").Value = 'クラウド(AWS/Azure) 運用管理 研修の演習サポート講師業務【経験不問】(再掲)'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5430951'
$ws.Range("G11").Value = 38
$ws.Range("H11").Value = '◇管理'
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5430951') | Out-Null

# Row 12
$ws.Range("A12").Value = '2025-11-11 01:21:02'
$ws.Range("B12").Value = '進行管理およびチームディレクションを担当'
$ws.Range("C12").Value = 'システム開発'
$ws.Range("D12").Value = '~ 5,000 円 / 固定'
$ws.Range("E12").Value = '期限情報なし'
$ws.Range("F12").Value = 'https://www.lancers.jp/work/detail/5418064'
$ws.Range("G12").Value = 30
$ws.Range("H12").Value = '◇管理'
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5418064') | Out-Null

# Row 13
$ws.Range("A13").Value = '2025-11-11 01:21:02'
$ws.Range("B13").Value = '初回 【高単価×長期案件あり】フリーランスエンジニア募集|リモート可・週3〜OK'
$ws.Range("C13").Value = 'システム開発'
$ws.Range("D13").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E13").Value = '期限情報なし'
$ws.Range("F13").Value = 'https://www.lancers.jp/work/detail/5431085'
$ws.Range("G13").Value = 25
$ws.Hyperlinks.Add($ws.Range("F13"), 'https://www.lancers.jp/work/detail/5431085') | Out-Null

# Row 14
$ws.Range("A14").Value = '2025-11-11 01:21:02'
$ws.Range("B14").Value = '【若手歓迎×リモートOK】SRE/インフラエンジニア(Google Cloud/長期・金融系案件)'
$ws.Range("C14").Value = 'システム開発'
$ws.Range("D14").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E14").Value = '期限情報なし'
$ws.Range("F14").Value = 'https://www.lancers.jp/work/detail/5428756'
$ws.Range("G14").Value = 25
$ws.Hyperlinks.Add($ws.Range("F14"), 'https://www.lancers.jp/work/detail/5428756') | Out-Null

# Row 15
$ws.Range("A15").Value = '2025-11-11 01:21:02'
$ws.Range("B15").Value = '【リーダー募集×リモートOK】SRE/インフラエンジニア(Google Cloud/長期金融系案件)'
$ws.Range("C15").Value = 'システム開発'
$ws.Range("D15").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E15").Value = '期限情報なし'
$ws.Range("F15").Value = 'https://www.lancers.jp/work/detail/5428755'
$ws.Range("G15").Value = 25
$ws.Hyperlinks.Add($ws.Range("F15"), 'https://www.lancers.jp/work/detail/5428755') | Out-Null

# Row 16
$ws.Range("A16").Value = '2025-11-11 01:21:02'
$ws.Range("B16").Value = 'AWS環境からAWS環境ヘの新規構築'
$ws.Range("C16").Value = 'システム開発'
$ws.Range("D16").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E16").Value = '期限情報なし'
$ws.Range("F16").Value = 'https://www.lancers.jp/work/detail/5431069'
$ws.Range("G16").Value = 18
$ws.Hyperlinks.Add($ws.Range("F16"), 'https://www.lancers.jp/work/detail/5431069') | Out-Null

# Row 17
$ws.Range("A17").Value = '2025-11-11 01:21:02'
$ws.Range("B17").Value = 'AMAZON運用代行 無在庫欧米輸入経験者のみ 物販知識がお有りの方'
$ws.Range("C17").Value = 'システム開発'
$ws.Range("D17").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E17").Value = '期限情報なし'
$ws.Range("F17").Value = 'https://www.lancers.jp/work/detail/5431036'
$ws.Range("G17").Value = 13
$ws.Hyperlinks.Add($ws.Range("F17"), 'https://www.lancers.jp/work/detail/5431036') | Out-Null

# Row 18
$ws.Range("A18").Value = '2025-11-11 01:21:02'
$ws.Range("B18").Value = 'EAの作成'
$ws.Range("C18").Value = 'システム開発'
$ws.Range("D18").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E18").Value = '期限情報なし'
$ws.Range("F18").Value = 'https://www.lancers.jp/work/detail/5431276'
$ws.Range("G18").Value = 10
$ws.Hyperlinks.Add($ws.Range("F18"), 'https://www.lancers.jp/work/detail/5431276') | Out-Null

# Row 19
$ws.Range("A19").Value = '2025-11-11 01:21:02'
$ws.Range("B19").Value = '【データ加工のプロ募集】施設情報データの修正・整備依頼'
$ws.Range("C19").Value = 'システム開発'
$ws.Range("D19").Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E19").Value = '期限情報なし'
$ws.Range("F19").Value = 'https://www.lancers.jp/work/detail/5417622'
$ws.Range("G19").Value = 10
$ws.Hyperlinks.Add($ws.Range("F19"), 'https://www.lancers.jp/work/detail/5417622') | Out-Null
